$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8527436852455139
$ws.Range("B1").Value = 2.850529432296753
$ws.Range("C1").Value = 3.158830642700195
$ws.Range("D1").Value = 1.784195184707642
$ws.Range("E1").Value = 1.368252754211426
